$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1010.9091
$ws.Range("I19").Value = 1014.1667
$ws.Range("J19").Value = 1007
$ws.Range("K19").Value = 1014.1667
$ws.Range("L19").Value = 1007
$ws.Range("M19").Value = -839.1667
$ws.Range("N19").Value = -1357
$ws.Range("H95").Value = 17812
$ws.Range("J95").Value = 17812
$ws.Range("L95").Value = 17812
$ws.Range("N95").Value = -23304
$ws.Range("H125").Value = 873.8946999999999
$ws.Range("I125").Value = 719.2222
$ws.Range("J125").Value = 1013.1
$ws.Range("K125").Value = 6472.999800000001
$ws.Range("L125").Value = 9117.9
$ws.Range("M125").Value = -4012.999800000001
$ws.Range("N125").Value = -14037.9
$ws.Range("H133").Value = 93840.664
$ws.Range("J133").Value = 93840.664
$ws.Range("L133").Value = 93840.664
$ws.Range("N133").Value = -103960.664
$ws.Range("H134").Value = 80525.11
$ws.Range("J134").Value = 91135.86
$ws.Range("L134").Value = 91135.86
$ws.Range("N134").Value = -101275.86
$ws.Range("H137").Value = 494571.38
$ws.Range("I137").Value = 1817.5
$ws.Range("J137").Value = 1300895.9
$ws.Range("K137").Value = 5452.5
$ws.Range("L137").Value = 3902687.7
$ws.Range("M137").Value = -2902.5
$ws.Range("N137").Value = -3907787.7
$ws.Range("H138").Value = 2652.3872
$ws.Range("J138").Value = 3192.0833
$ws.Range("L138").Value = 9576.249899999999
$ws.Range("N138").Value = -19856.2499

# Sheet 2: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H63").Value = 3301.4
$ws.Range("I63").Value = 2167
$ws.Range("K63").Value = 2167
$ws.Range("M63").Value = -1481
$ws.Range("H66").Value = 3301.4
$ws.Range("I66").Value = 2167
$ws.Range("K66").Value = 10835
$ws.Range("M66").Value = -7403
$ws.Range("H102").Value = 45308.89
$ws.Range("J102").Value = 61056
$ws.Range("L102").Value = 61056
$ws.Range("N102").Value = -64300
$ws.Range("H132").Value = 2353.4211
$ws.Range("I132").Value = 1866.138
$ws.Range("J132").Value = 3923.5557
$ws.Range("K132").Value = 5598.414
$ws.Range("L132").Value = 11770.6671
$ws.Range("M132").Value = -3068.414
$ws.Range("N132").Value = -16830.6671
$ws.Range("H139").Value = 105000
$ws.Range("J139").Value = 105000
$ws.Range("L139").Value = 105000
$ws.Range("N139").Value = -115280

# Sheet 3: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 110656.26
$ws.Range("I20").Value = 178266.69
$ws.Range("J20").Value = 3606.4167
$ws.Range("K20").Value = 178266.69
$ws.Range("L20").Value = 3606.4167
$ws.Range("M20").Value = -178019.69
$ws.Range("N20").Value = -4100.4167
$ws.Range("H86").Value = 1805.3182
$ws.Range("I86").Value = 1138.875
$ws.Range("J86").Value = 3582.5
$ws.Range("K86").Value = 1138.875
$ws.Range("L86").Value = 3582.5
$ws.Range("M86").Value = -15.875
$ws.Range("N86").Value = -5828.5
$ws.Range("H89").Value = 1805.3182
$ws.Range("I89").Value = 1138.875
$ws.Range("J89").Value = 3582.5
$ws.Range("K89").Value = 5694.375
$ws.Range("L89").Value = 17912.5
$ws.Range("M89").Value = -78.375
$ws.Range("N89").Value = -29144.5
$ws.Range("H105").Value = 2855.5833
$ws.Range("I105").Value = 2100.6
$ws.Range("K105").Value = 2100.6
$ws.Range("M105").Value = -353.5999999999999

# Sheet 4: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 464.73334
$ws.Range("I22").Value = 464.73334
$ws.Range("K22").Value = 464.73334
$ws.Range("M22").Value = -114.73334
$ws.Range("H31").Value = 4371.143
$ws.Range("I31").Value = 2358
$ws.Range("J31").Value = 5994.645
$ws.Range("K31").Value = 2358
$ws.Range("L31").Value = 5994.645
$ws.Range("M31").Value = -2063
$ws.Range("N31").Value = -6584.645
$ws.Range("H34").Value = 4371.143
$ws.Range("I34").Value = 2358
$ws.Range("J34").Value = 5994.645
$ws.Range("K34").Value = 2358
$ws.Range("L34").Value = 5994.645
$ws.Range("M34").Value = -2156
$ws.Range("N34").Value = -6398.645
$ws.Range("H58").Value = 2227.4546
$ws.Range("I58").Value = 1966
$ws.Range("K58").Value = 1966
$ws.Range("M58").Value = -1763
$ws.Range("H99").Value = 1955381.9
$ws.Range("I99").Value = 1723.75
$ws.Range("K99").Value = 1723.75
$ws.Range("M99").Value = -225.75
$ws.Range("H126").Value = 1955381.9
$ws.Range("I126").Value = 1723.75
$ws.Range("K126").Value = 5171.25
$ws.Range("M126").Value = -2701.25
$ws.Range("H132").Value = 1421.75
$ws.Range("I132").Value = 996.2222
$ws.Range("K132").Value = 2988.6666
$ws.Range("M132").Value = -458.6666
$ws.Range("H134").Value = 44399.168
$ws.Range("I134").Value = 2685.611
$ws.Range("J134").Value = 169539.83
$ws.Range("K134").Value = 8056.833
$ws.Range("L134").Value = 508619.49
$ws.Range("M134").Value = -5521.833
$ws.Range("N134").Value = -513689.49
$ws.Range("H136").Value = 2227.4546
$ws.Range("I136").Value = 1966
$ws.Range("K136").Value = 5898
$ws.Range("M136").Value = -3348
$ws.Range("H138").Value = 55282.5
$ws.Range("J138").Value = 54894.285
$ws.Range("L138").Value = 54894.285
$ws.Range("N138").Value = -65174.285
$ws.Range("H141").Value = 188117.66
$ws.Range("J141").Value = 206163.67
$ws.Range("L141").Value = 206163.67
$ws.Range("N141").Value = -216523.67

# Sheet 5: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 1499.6666
$ws.Range("I58").Value = 499
$ws.Range("K58").Value = 1497
$ws.Range("M58").Value = -1369
$ws.Range("H98").Value = 400
$ws.Range("I98").Value = 300
$ws.Range("J98").Value = 433.33334
$ws.Range("K98").Value = 900
$ws.Range("L98").Value = 1300.00002
$ws.Range("N98").Value = -4296.000019999999
$ws.Range("M98").Value = 598
$ws.Range("H121").Value = 556983.6
$ws.Range("I121").Value = 516.2857
$ws.Range("K121").Value = 1548.8571
$ws.Range("M121").Value = -238.8571000000002
$ws.Range("H126").Value = 5025.375
$ws.Range("I126").Value = 5025.375
$ws.Range("K126").Value = 15076.125
$ws.Range("M126").Value = -10136.125

# Sheet 6: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 38298.258
$ws.Range("I70").Value = 49187
$ws.Range("K70").Value = 49187
$ws.Range("M70").Value = -48917
$ws.Range("H73").Value = 38298.258
$ws.Range("I73").Value = 49187
$ws.Range("K73").Value = 49187
$ws.Range("M73").Value = -48251
$ws.Range("H116").Value = 50767.69
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 50767.69
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 50767.69
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -59945.69
$ws.Range("H122").Value = 131439.27
$ws.Range("I122").Value = 142021.62
$ws.Range("K122").Value = 426064.86
$ws.Range("M122").Value = -423614.86
$ws.Range("H132").Value = 2325.9644
$ws.Range("I132").Value = 2128.5
$ws.Range("K132").Value = 6385.5
$ws.Range("M132").Value = -3855.5
$ws.Range("H135").Value = 76618.625
$ws.Range("J135").Value = 76618.625
$ws.Range("L135").Value = 76618.625
$ws.Range("N135").Value = -86758.625

# Sheet 7: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1346664.4
$ws.Range("I2").Value = 1696.4
$ws.Range("K2").Value = 1696.4
$ws.Range("M2").Value = -1584.4
$ws.Range("H68").Value = 701800.3
$ws.Range("I68").Value = 701800.3
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 701800.3
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -701051.3
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 701800.3
$ws.Range("I71").Value = 701800.3
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 3509001.5
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -3505257.5
$ws.Range("N71").ClearContents()
$ws.Range("H122").Value = 40007068
$ws.Range("J122").Value = 66672612
$ws.Range("L122").Value = 200017836
$ws.Range("N122").Value = -200022736
$ws.Range("H132").Value = 2912.6785
$ws.Range("I132").Value = 2414.75
$ws.Range("J132").Value = 3576.5833
$ws.Range("K132").Value = 7244.25
$ws.Range("L132").Value = 10729.7499
$ws.Range("M132").Value = -4714.25
$ws.Range("N132").Value = -15789.7499
$ws.Range("H136").Value = 6486
$ws.Range("I136").Value = 6644.1113
$ws.Range("K136").Value = 19932.3339
$ws.Range("M136").Value = -17382.3339

# Sheet 8: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1177007.2
$ws.Range("I132").Value = 1831.7587
$ws.Range("J132").Value = 5437018.5
$ws.Range("K132").Value = 5495.2761
$ws.Range("L132").Value = 16311055.5
$ws.Range("M132").Value = -2965.2761
$ws.Range("N132").Value = -16316115.5
